$wb = $excel.ActiveWorkbook

# New cell style used for the pasted-in data tables (font size 12, based on "Normalny").
$dataStyle = $wb.Styles.Add("Normalny 2")
$dataStyle.Font.Size = 12

# ---------------------------------------------------------------------------
# Arkusz7: Rok / Ilość
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$s7 = $wb.Worksheets.Add($null, $lastSheet)
$s7.Name = "Arkusz7"

$s7.Range("A1").Value = "Rok"
$s7.Range("B1").Value = "Ilość"

$s7.Range("A2").Value = 2019
$s7.Range("B2").Value = 57
$s7.Range("A3").Value = 2020
$s7.Range("B3").Value = 201
$s7.Range("A4").Value = 2021
$s7.Range("B4").Value = 220

$s7.Range("A2:B4").Style = "Normalny 2"
$s7.Range("A2:B4").RowHeight = 15.75

[void]$s7.Range("B1").Select()

# ---------------------------------------------------------------------------
# Arkusz8: Data / zakup
# ---------------------------------------------------------------------------
$s8 = $wb.Worksheets.Add($null, $s7)
$s8.Name = "Arkusz8"

$s8.Range("A1").Value = "Data"
$s8.Range("B1").Value = "zakup"

$s8.Range("A2").Value = 2017
$s8.Range("B2").Value = 12
$s8.Range("A3").Value = 2018
$s8.Range("B3").Value = 17
$s8.Range("A4").Value = 2019
$s8.Range("B4").Value = 28
$s8.Range("A5").Value = 2020
$s8.Range("B5").Value = 29.4
$s8.Range("A6").Value = 2021
$s8.Range("B6").Value = 32.3

$s8.Range("A2:B6").Style = "Normalny 2"
$s8.Range("A2:B6").RowHeight = 15.75

[void]$s8.Range("B1").Select()

# ---------------------------------------------------------------------------
# Arkusz9: Data / Zakup
# ---------------------------------------------------------------------------
$s9 = $wb.Worksheets.Add($null, $s8)
$s9.Name = "Arkusz9"

$s9.Range("A1").Value = "Data"
$s9.Range("B1").Value = "Zakup"

$s9.Range("A2").Value = 2017
$s9.Range("B2").Value = 44
$s9.Range("A3").Value = 2018
$s9.Range("B3").Value = 42
$s9.Range("A4").Value = 2019
$s9.Range("B4").Value = 31
$s9.Range("A5").Value = 2020
$s9.Range("B5").Value = 26
$s9.Range("A6").Value = 2021
$s9.Range("B6").Value = 22

$s9.Range("A2:B6").Style = "Normalny 2"
$s9.Range("A2:B6").RowHeight = 15.75

[void]$s9.Range("C2").Select()

# ---------------------------------------------------------------------------
# Arkusz10: Data / zasieg
# ---------------------------------------------------------------------------
$s10 = $wb.Worksheets.Add($null, $s9)
$s10.Name = "Arkusz10"

$s10.Range("A1").Value = "Data"
$s10.Range("B1").Value = "zasieg"

$s10.Range("A2").Value = 2017
$s10.Range("B2").Value = 228
$s10.Range("A3").Value = 2018
$s10.Range("B3").Value = 273
$s10.Range("A4").Value = 2019
$s10.Range("B4").Value = 309
$s10.Range("A5").Value = 2020
$s10.Range("B5").Value = 320
$s10.Range("A6").Value = 2021
$s10.Range("B6").Value = 390

[void]$s10.Range("C3").Select()
